$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "Designation" header (column K) to "Speciality"
$kCell = $ws.Range("K1")
$kCell.Value = "Speciality"

# Give the K1 header cell its own explicit alignment formatting (same font as
# before) instead of sharing the plain, no-alignment style used by the other
# header cells.
$kCell.WrapText = $false
